$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title shape: "A" + " " + "slide" (3 runs) -> "A slide" (1 run) ---
# Re-assigning the exact same text is treated as a no-op by the engine's
# run-diffing, so the old run split would survive untouched. Flipping the
# text to something unrelated first forces a genuine replace (no shared
# prefix/suffix with the old runs), consolidating the whole paragraph into
# a single fresh run; setting the final text immediately afterwards then
# just edits that single run in place, so it keeps the plain <a:rPr/> with
# no extra attributes such as lang.
$titleShape = $s.Shapes.Item(1)
$titleTextRange = $titleShape.TextFrame.TextRange
$titleTextRange.Text = "~"
$titleTextRange.Text = "A slide"
